$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.956.47'
$ws.Range("E2").Value = '  -2.25%  '

$ws.Range("D3").Value = '1.867.33'
$ws.Range("E3").Value = '  -2.53%  '

$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = '312.21'
$ws.Range("E5").Value = '  -1.30%  '

$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.02%  '

$ws.Range("D7").Value = '0.4984'
$ws.Range("E7").Value = '  -3.40%  '

$ws.Range("D8").Value = '0.3812'
$ws.Range("E8").Value = '  -4.73%  '

$ws.Range("D9").Value = '0.08946'
$ws.Range("E9").Value = '  -8.75%  '

$ws.Range("D10").Value = '1.119'
$ws.Range("E10").Value = '  -2.77%  '

$ws.Range("D11").Value = '41.43'
$ws.Range("E11").Value = '  -1.95%  '

$ws.Range("D12").Value = '6.329'
$ws.Range("E12").Value = '  -2.89%  '

$ws.Range("D13").Value = '20.67'
$ws.Range("E13").Value = '  -2.69%  '

$ws.Range("D14").Value = '1.870.26'
$ws.Range("E14").Value = '  -1.92%  '

$ws.Range("D15").Value = '7.236'
$ws.Range("E15").Value = '  -3.20%  '

$ws.Range("D16").Value = '1.002'
$ws.Range("E16").Value = '  +0.03%  '

$ws.Range("D17").Value = "'0.00001100"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -3.22%  '

$ws.Range("D18").Value = '90.79'
$ws.Range("E18").Value = '  -4.13%  '

$ws.Range("D19").Value = '0.06654'
$ws.Range("E19").Value = '  -0.04%  '

$ws.Range("D20").Value = '17.97'
$ws.Range("E20").Value = '  -1.61%  '

$ws.Range("D21").Value = "'1.000"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.03%  '

$ws.Range("D22").Value = '6.115'
$ws.Range("E22").Value = '  -3.08%  '

$ws.Range("D23").Value = '27.980.64'
$ws.Range("E23").Value = '  -2.34%  '

$ws.Range("D24").Value = '11.58'
$ws.Range("E24").Value = '  +0.54%  '

$ws.Range("D25").Value = '2.286'
$ws.Range("E25").Value = '  -1.17%  '

$ws.Range("D26").Value = '2.084.52'
$ws.Range("E26").Value = '  -2.65%  '

$ws.Range("D27").Value = '2.508'
$ws.Range("E27").Value = '  -6.48%  '

$ws.Range("D28").Value = '158.47'
$ws.Range("E28").Value = '  +0.43%  '

$ws.Range("D29").Value = '20.68'
$ws.Range("E29").Value = '  -2.90%  '

$ws.Range("D30").Value = '126.31'
$ws.Range("E30").Value = '  -2.27%  '

$ws.Range("D31").Value = '0.1058'
$ws.Range("E31").Value = '  -1.40%  '

$ws.Range("D32").Value = '1.055'
$ws.Range("E32").Value = '  -5.49%  '

$ws.Range("D33").Value = '5.574'
$ws.Range("E33").Value = '  -2.61%  '

$ws.Range("D34").Value = "'3.580"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -1.40%  '

$ws.Range("D35").Value = '9.424'
$ws.Range("E35").Value = '  -4.89%  '

$ws.Range("D36").Value = '0.06542'
$ws.Range("E36").Value = '  -3.68%  '

$ws.Range("D37").Value = "'0.02400"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.72%  '

$ws.Range("D38").Value = "'1.300"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +9.53%  '

$ws.Range("D39").Value = '0.2191'
$ws.Range("E39").Value = '  -1.90%  '

$ws.Range("D40").Value = "'1.200"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -5.56%  '

$ws.Range("D41").Value = '11.68'
$ws.Range("E41").Value = '  -1.39%  '

$ws.Range("D42").Value = '0.6365'
$ws.Range("E42").Value = '  -2.07%  '

$ws.Range("E43").Value = '  -3.64%  '

$ws.Range("D44").Value = "'1.000"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.01%  '

$ws.Range("D45").Value = '13.18'
$ws.Range("E45").Value = '  -3.92%  '

$ws.Range("D46").Value = '0.6004'
$ws.Range("E46").Value = '  -1.74%  '

$ws.Range("D47").Value = '1.278'
$ws.Range("E47").Value = '  -0.78%  '

$ws.Range("D48").Value = '3.674'
$ws.Range("E48").Value = '  -2.97%  '

$ws.Range("D49").Value = '1.224'
$ws.Range("E49").Value = '  +1.52%  '

$ws.Range("D50").Value = '1.985'
$ws.Range("E50").Value = '  -3.99%  '

$ws.Range("D51").Value = '121.43'
$ws.Range("E51").Value = '  -2.86%  '
